# Fill in the next working-time entry (row 35) on the single worksheet.
# Mirrors the pattern of the existing rows (e.g. row 34): a date in column E
# (formatted like the others), hours worked in F, unit "Stunden" in G,
# category "Programmieren" in H, and a description of the task in I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy E34's formatting (date number format) onto E35, then set its value -
# this reuses the existing date style instead of minting a new number format.
$ws.Range("E34").Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("E35").Value = 43691          # 2019-08-14
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = "Stunden"
$ws.Range("H35").Value = "Programmieren"
$ws.Range("I35").Value = "Lösen von Dependency Fehlern"

# Match the selection state recorded in the workbook after the edit.
$ws.Range("I35").Select() | Out-Null
